# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from 45543 (2024-09-08) to 45544 (2024-09-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45543) {
        $cell.Value = 45544
    }
}
